$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Table 1 (top table, Solubility) ---

# Header row 2: extend with three more repeats of r2 score/mse/rmse/mae
# (one block per additional model variant added below).
$ws.Range("G2").Value2 = "r2 score"
$ws.Range("H2").Value2 = "mse"
$ws.Range("I2").Value2 = "rmse"
$ws.Range("J2").Value2 = "mae"
$ws.Range("K2").Value2 = "r2 score"
$ws.Range("L2").Value2 = "mse"
$ws.Range("M2").Value2 = "rmse"
$ws.Range("N2").Value2 = "mae"

# Drop the old "Kfold"/"No" annotation cells (G3, G4) - no longer used.
$ws.Range("G3").ClearContents()
$ws.Range("G4").ClearContents()

# Row 5 keeps "Mordred RF" but now carries full metrics across three more
# model variants (columns C..N).
$ws.Range("A5").Value2 = "Mordred RF"
$ws.Range("C5").Value2 = 0.821875364512533
$ws.Range("D5").Value2 = 0.93363946285174
$ws.Range("E5").Value2 = 0.965841640884103
$ws.Range("F5").Value2 = 0.609196229486432
$ws.Range("G5").Value2 = 0.821720294387933
$ws.Range("H5").Value2 = 0.934424643127384
$ws.Range("I5").Value2 = 0.96626140041123
$ws.Range("J5").Value2 = 0.609728853729244
$ws.Range("K5").Value2 = 0.821728229549132
$ws.Range("L5").Value2 = 0.934405591064197
$ws.Range("M5").Value2 = 0.966242173125073
$ws.Range("N5").Value2 = 0.609520694448872

# Replace the old "Mordred NN" row (row 6) and add new rows for the
# additional featurization/model comparisons (rows 7-19 were blank).
$ws.Range("A6").Value2 = "RDKIT RF"
$ws.Range("A7").Value2 = "Mol2Vec RF"
$ws.Range("A8").Value2 = "ecfp RF"
$ws.Range("A9").Value2 = "pubchem RF"
$ws.Range("A10").Value2 = "maccs key"
$ws.Range("A11").Value2 = "spectrophore"

# --- Table 2 (lower table, Cocrystal) ---

# Row 23 ("Resnet18"): clear the stale B/C/E values, update D to the new
# best accuracy figure.
$ws.Range("B23").ClearContents()
$ws.Range("C23").ClearContents()
$ws.Range("D23").Value2 = 0.753999984264374
$ws.Range("E23").ClearContents()

# --- View state ---
$ws.Range("D11").Select()
